# Remove stale/incorrect w:link back-references from a handful of styles.
#
# The source document carries a number of character/paragraph style pairs
# whose <w:link> cross-reference points at a style id/name that no longer
# exists in the style sheet (e.g. "EntteCar" links to "En-tte", which isn't
# a real style). Those dangling links are simply wrong and get cleared;
# the still-valid linked pairs (Heading1/2/3 <-> Titre1/2/3Car, and the
# Footer/BalloonText/Title paragraph styles that correctly point at their
# character-style counterparts) are left untouched.

$d = $word.ActiveDocument

$staleLinkedStyles = @(
    "EntteCar",
    "PieddepageCar",
    "TextedebullesCar",
    "Titre1Car",
    "Titre2Car",
    "Titre3Car",
    "TitreCar",
    "Header"
)

foreach ($styleName in $staleLinkedStyles) {
    $style = $d.Styles($styleName)
    $style.LinkStyle = $null
}
